# Auto-generated edit script: update currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across multiple sheets per scheduled Universalis price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4447343.5
$ws.Range("I33").Value = 5307380
$ws.Range("K33").Value = 5307380
$ws.Range("M33").Value = -5307151

$ws.Range("H97").Value = 3822.5
$ws.Range("J97").Value = 4478
$ws.Range("L97").Value = 13434
$ws.Range("N97").Value = -14426

$ws.Range("H98").Value = 4124.25
$ws.Range("I98").Value = 4124.25
$ws.Range("K98").Value = 4124.25
$ws.Range("M98").Value = -2626.25

$ws.Range("H101").Value = 1644.9231
$ws.Range("I101").Value = 1709.8889
$ws.Range("J101").Value = 1498.75
$ws.Range("K101").Value = 5129.6667
$ws.Range("L101").Value = 4496.25
$ws.Range("M101").Value = -3507.6667
$ws.Range("N101").Value = -7740.25

$ws.Range("H122").Value = 4124.25
$ws.Range("I122").Value = 4124.25
$ws.Range("K122").Value = 12372.75
$ws.Range("M122").Value = -9922.75

$ws.Range("H135").Value = 1649.75
$ws.Range("I135").Value = 800
$ws.Range("K135").Value = 7200
$ws.Range("M135").Value = -4665

$ws.Range("H137").Value = 2332522.8
$ws.Range("I137").Value = 4924.6294
$ws.Range("J137").Value = 6260344.5
$ws.Range("K137").Value = 14773.8882
$ws.Range("L137").Value = 18781033.5
$ws.Range("M137").Value = -12223.8882
$ws.Range("N137").Value = -18786133.5

$ws.Range("H138").Value = 5099.8125
$ws.Range("I138").Value = 5973.143
$ws.Range("K138").Value = 17919.429
$ws.Range("M138").Value = -12779.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10823.667
$ws.Range("I28").Value = 10823.667
$ws.Range("K28").Value = 10823.667
$ws.Range("M28").Value = -10631.667

$ws.Range("H32").Value = 307410.88
$ws.Range("I32").Value = 401841.97
$ws.Range("K32").Value = 401841.97
$ws.Range("M32").Value = -401554.97

$ws.Range("H99").Value = 10823.667
$ws.Range("I99").Value = 10823.667
$ws.Range("K99").Value = 10823.667
$ws.Range("M99").Value = -7828.666999999999

$ws.Range("H132").Value = 3020.2666
$ws.Range("I132").Value = 2576.0667
$ws.Range("K132").Value = 7728.2001
$ws.Range("M132").Value = -5198.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 983.375
$ws.Range("I22").Value = 688.1429000000001
$ws.Range("K22").Value = 688.1429000000001
$ws.Range("M22").Value = -515.1429000000001

$ws.Range("H103").Value = 35542
$ws.Range("J103").Value = 35542
$ws.Range("L103").Value = 35542
$ws.Range("N103").Value = -37886

$ws.Range("H134").Value = 21952884
$ws.Range("I134").Value = 1624.3462
$ws.Range("J134").Value = 60001732
$ws.Range("K134").Value = 4873.0386
$ws.Range("L134").Value = 180005196
$ws.Range("M134").Value = -2338.0386
$ws.Range("N134").Value = -180010266

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2523.9592
$ws.Range("I31").Value = 1980.4878
$ws.Range("K31").Value = 1980.4878
$ws.Range("M31").Value = -1685.4878

$ws.Range("H34").Value = 2523.9592
$ws.Range("I34").Value = 1980.4878
$ws.Range("K34").Value = 1980.4878
$ws.Range("M34").Value = -1778.4878

$ws.Range("H51").Value = 69999
$ws.Range("J51").Value = 69999
$ws.Range("L51").Value = 69999
$ws.Range("N51").Value = -71471

$ws.Range("H61").Value = 69999
$ws.Range("J61").Value = 69999
$ws.Range("L61").Value = 69999
$ws.Range("N61").Value = -70695

$ws.Range("H86").Value = 21408.908
$ws.Range("I86").Value = 11400
$ws.Range("J86").Value = 29749.666
$ws.Range("K86").Value = 11400
$ws.Range("L86").Value = 29749.666
$ws.Range("M86").Value = -10277
$ws.Range("N86").Value = -31995.666

$ws.Range("H89").Value = 21408.908
$ws.Range("I89").Value = 11400
$ws.Range("J89").Value = 29749.666
$ws.Range("K89").Value = 57000
$ws.Range("L89").Value = 148748.33
$ws.Range("M89").Value = -51384
$ws.Range("N89").Value = -159980.33

$ws.Range("H106").Value = 54932.2
$ws.Range("J106").Value = 54932.2
$ws.Range("L106").Value = 54932.2
$ws.Range("N106").Value = -57456.2

$ws.Range("H134").Value = 1794
$ws.Range("I134").Value = 1631
$ws.Range("J134").Value = 2120
$ws.Range("K134").Value = 4893
$ws.Range("L134").Value = 6360
$ws.Range("M134").Value = -2358
$ws.Range("N134").Value = -11430

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14455.5625
$ws.Range("I3").Value = 12642.833
$ws.Range("K3").Value = 37928.499
$ws.Range("M3").Value = -37816.499

$ws.Range("H4").Value = 43024.32
$ws.Range("I4").Value = 487.36365
$ws.Range("J4").Value = 666899.7
$ws.Range("K4").Value = 1462.09095
$ws.Range("L4").Value = 2000699.1
$ws.Range("M4").Value = -1350.09095
$ws.Range("N4").Value = -2000923.1

$ws.Range("H5").Value = 889.3333
$ws.Range("J5").Value = 1166.6666
$ws.Range("L5").Value = 3499.9998
$ws.Range("N5").Value = -3723.9998

$ws.Range("H68").Value = 3083.6
$ws.Range("J68").Value = 3497.5
$ws.Range("L68").Value = 10492.5
$ws.Range("N68").Value = -12114.5

$ws.Range("H71").Value = 3083.6
$ws.Range("J71").Value = 3497.5
$ws.Range("L71").Value = 31477.5
$ws.Range("N71").Value = -39589.5

$ws.Range("H122").Value = 3705068.8
$ws.Range("J122").Value = 3538
$ws.Range("L122").Value = 31842
$ws.Range("N122").Value = -36742

$ws.Range("H134").Value = 1994.909
$ws.Range("J134").Value = 3750
$ws.Range("L134").Value = 11250
$ws.Range("N134").Value = -21390

$ws.Range("H135").Value = 889.3333
$ws.Range("J135").Value = 1166.6666
$ws.Range("L135").Value = 10499.9994
$ws.Range("N135").Value = -15569.9994

$ws.Range("H136").Value = 8710
$ws.Range("I136").Value = 3719.5454
$ws.Range("K136").Value = 11158.6362
$ws.Range("M136").Value = -6058.636200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4254.6
$ws.Range("I41").Value = 4254.6
$ws.Range("K41").Value = 4254.6
$ws.Range("M41").Value = -3899.6

$ws.Range("H99").Value = 9620.817999999999
$ws.Range("I99").Value = 7083
$ws.Range("K99").Value = 7083
$ws.Range("M99").Value = -4837

$ws.Range("H102").Value = 41668520
$ws.Range("I102").Value = 50001824
$ws.Range("J102").Value = 1999.5
$ws.Range("K102").Value = 50001824
$ws.Range("L102").Value = 1999.5
$ws.Range("M102").Value = -50000202
$ws.Range("N102").Value = -5243.5

$ws.Range("H105").Value = 49149.168
$ws.Range("J105").Value = 49149.168
$ws.Range("L105").Value = 49149.168
$ws.Range("N105").Value = -56137.168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 45000
$ws.Range("J3").Value = 45000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45224

$ws.Range("H15").Value = 45000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45340

$ws.Range("H20").Value = 367000000
$ws.Range("I20").Value = 999998
$ws.Range("K20").Value = 999998
$ws.Range("M20").Value = -999772

$ws.Range("H22").Value = 4270.95
$ws.Range("J22").Value = 5041.9375
$ws.Range("L22").Value = 5041.9375
$ws.Range("N22").Value = -5631.9375

$ws.Range("H27").Value = 4270.95
$ws.Range("J27").Value = 5041.9375
$ws.Range("L27").Value = 5041.9375
$ws.Range("N27").Value = -5255.9375

$ws.Range("H93").Value = 1688
$ws.Range("I93").Value = 1688
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1688
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -440
$ws.Range("N93").ClearContents()

$ws.Range("H99").Value = 46552.453
$ws.Range("I99").Value = 31207.7
$ws.Range("K99").Value = 31207.7
$ws.Range("M99").Value = -28212.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 9500
$ws.Range("J58").Value = 8000
$ws.Range("L58").Value = 8000
$ws.Range("N58").Value = -8616

$ws.Range("H61").Value = 8823.25
$ws.Range("I61").Value = 8823.25
$ws.Range("K61").Value = 8823.25
$ws.Range("M61").Value = -8531.25

$ws.Range("H76").Value = 109829.836
$ws.Range("J76").Value = 109829.836
$ws.Range("L76").Value = 109829.836
$ws.Range("N76").Value = -110459.836

$ws.Range("H79").Value = 109829.836
$ws.Range("J79").Value = 109829.836
$ws.Range("L79").Value = 109829.836
$ws.Range("N79").Value = -112013.836

$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

$ws.Range("H139").Value = 86500
$ws.Range("J139").Value = 86500
$ws.Range("L139").Value = 86500
$ws.Range("N139").Value = -96780
